# Add a new "PRIDE_PROTEOMICS" ER (Existing Resource) sheet as the last
# sheet of the workbook, fill it with the PRIDE proteomics ER table
# content, and make it the active/selected sheet - mirrors the commit
# "Fill empty ER sheets into every template (except Imaging)".

$wb = $excel.ActiveWorkbook

# Add the new worksheet directly after the current last sheet so it ends
# up as the 3rd / final tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PRIDE_PROTEOMICS"
$ws = $newSheet

# Pre-fill the whole used range (A1:L13) with an explicit empty-text value
# (quote-prefixed empty string) so every cell - even the "blank" ones - is
# a real empty string cell, matching the source template's cells, then
# strip the quote-prefix formatting that trick adds.
$ws.Range("A1:L13").Value = "'"
$ws.Cells.ClearFormats()

# Header row
$ws.Range("B1").Value = "TermSourceRef"
$ws.Range("C1").Value = "Ontology"
$ws.Range("D1").Value = "TAN"
$ws.Range("E1").Value = "Content type (validation)"
$ws.Range("F1").Value = "Notes during templating"
$ws.Range("G1").Value = "Target term"
$ws.Range("H1").Value = "Instruction"
$ws.Range("I1").Value = "Requirement (m/o/n)"
$ws.Range("J1").Value = "Value (cv/s/d)"
$ws.Range("K1").Value = "Additional information"
$ws.Range("L1").Value = "Review comments"

# Data rows - one ER term per row
$ws.Range("A2").Value = "Source Name"

$ws.Range("A3").Value = "Parameter [Experiment type]"
$ws.Range("B3").Value = "PRIDE:0000457"
$ws.Range("C3").Value = "PRIDE"
$ws.Range("D3").Value = "http://purl.obolibrary.org/obo/PRIDE_0000457"

$ws.Range("A4").Value = "Parameter [Quantification method]"
$ws.Range("B4").Value = "user-specific"
$ws.Range("C4").Value = "user-specific"
$ws.Range("D4").Value = "user-specific"

$ws.Range("A5").Value = "Parameter [technical replicate]"
$ws.Range("B5").Value = "MS:1001808"
$ws.Range("C5").Value = "MS"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/MS_1001808"

$ws.Range("A6").Value = "Parameter [Variable modification]"
$ws.Range("B6").Value = "user-specific"
$ws.Range("C6").Value = "user-specific"
$ws.Range("D6").Value = "user-specific"

$ws.Range("A7").Value = "Parameter [Fixed modification]"
$ws.Range("B7").Value = "user-specific"
$ws.Range("C7").Value = "user-specific"
$ws.Range("D7").Value = "user-specific"

$ws.Range("A8").Value = "Parameter [sample volume]"
$ws.Range("B8").Value = "MS:1000005"
$ws.Range("C8").Value = "MS"
$ws.Range("D8").Value = "http://purl.obolibrary.org/obo/MS_1000005"

$ws.Range("A9").Value = "Parameter [injection volume]"
$ws.Range("B9").Value = "user-specific"
$ws.Range("C9").Value = "user-specific"
$ws.Range("D9").Value = "user-specific"

$ws.Range("A10").Value = "Parameter [count unit]"
$ws.Range("B10").Value = "UO:0000189"
$ws.Range("C10").Value = "UO"
$ws.Range("D10").Value = "http://purl.obolibrary.org/obo/UO_0000189"

$ws.Range("A11").Value = "Parameter [instrument model]"
$ws.Range("B11").Value = "MS:1000031"
$ws.Range("C11").Value = "MS"
$ws.Range("D11").Value = "http://purl.obolibrary.org/obo/MS_1000031"

$ws.Range("A12").Value = "Parameter [duration]"
$ws.Range("B12").Value = "PATO:0001309"
$ws.Range("C12").Value = "PATO"
$ws.Range("D12").Value = "http://purl.obolibrary.org/obo/PATO_0001309"

$ws.Range("A13").Value = "Data File Name"

# Column widths (best-fit approximations of the original template's
# bestFit column widths).
$ws.Columns.Item(1).ColumnWidth = 32.022135416666664
$ws.Columns.Item(2).ColumnWidth = 13.736979166666666
$ws.Columns.Item(3).ColumnWidth = 11.451822916666666
$ws.Columns.Item(4).ColumnWidth = 42.022135416666664
$ws.Columns.Item(5).ColumnWidth = 22.736979166666668
$ws.Columns.Item(6).ColumnWidth = 22.022135416666668
$ws.Columns.Item(7).ColumnWidth = 10.451822916666666
$ws.Columns.Item(8).ColumnWidth = 9.736979166666666
$ws.Columns.Item(9).ColumnWidth = 19.592447916666668
$ws.Columns.Item(10).ColumnWidth = 12.736979166666666
$ws.Columns.Item(11).ColumnWidth = 20.736979166666668
$ws.Columns.Item(12).ColumnWidth = 16.592447916666668

# Select the whole sheet (matches the captured "select all" state) and
# make this new sheet the active tab of the workbook.
$ws.Cells.Select()
$ws.Activate()
